# Penalty Reward System (unfinished) -- shift the weekly forecast window
# forward by one week and overwrite MyForecast with the (partial/unfinished)
# penalty-reward numbers; mirror a few of the now-stale figures into the
# Summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Forecast Comparison": each row's Week_Start_Date (col B) rolls
# forward to the following week's date, and MyForecast (col D) is replaced
# with the new (mostly 0/1) value.
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$weekStartDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$myForecast = @(1, 1, 1, 0, 1, 1, 1, 0, 1, 1, 0, 1, 1, 1, 1, 1)

$firstRow = 2
for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $firstRow + $i

    # Force column B to stay text so the date string isn't reinterpreted as
    # a date serial number.
    $wsForecast.Cells.Item($row, 2).NumberFormat = "@"
    $wsForecast.Cells.Item($row, 2).Value = $weekStartDates[$i]

    $wsForecast.Cells.Item($row, 4).Value = $myForecast[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Summary": a handful of derived metrics that reference the shifted
# window / updated forecast values.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Cells.Item(2, 2).Value = "2023-01-15 to 2025-01-05"

$wsSummary.Cells.Item(4, 2).NumberFormat = "@"
$wsSummary.Cells.Item(4, 2).Value = "85"

$wsSummary.Cells.Item(6, 2).NumberFormat = "@"
$wsSummary.Cells.Item(6, 2).Value = "9"

$wsSummary.Cells.Item(7, 2).NumberFormat = "@"
$wsSummary.Cells.Item(7, 2).Value = "29"

$wsSummary.Cells.Item(8, 2).Value = "2300 units"

$wsSummary.Cells.Item(9, 2).NumberFormat = "@"
$wsSummary.Cells.Item(9, 2).Value = "11"

$wsSummary.Cells.Item(10, 2).NumberFormat = "@"
$wsSummary.Cells.Item(10, 2).Value = "5"

$wsSummary.Cells.Item(11, 2).NumberFormat = "@"
$wsSummary.Cells.Item(11, 2).Value = "2"

$wsSummary.Cells.Item(12, 2).NumberFormat = "@"
$wsSummary.Cells.Item(12, 2).Value = "1"

$wsSummary.Cells.Item(13, 2).NumberFormat = "@"
$wsSummary.Cells.Item(13, 2).Value = "2025-04-20"

$wsSummary.Cells.Item(14, 2).NumberFormat = "@"
$wsSummary.Cells.Item(14, 2).Value = "0"

$wsSummary.Cells.Item(15, 2).NumberFormat = "@"
$wsSummary.Cells.Item(15, 2).Value = "2025-02-02"
